$d = $word.ActiveDocument

# 1) Contacts loop body: {{c.name}} | {{c.phone}} | {{c.relation}}
#    -> {{INS $c.name}} | {{INS $c.phone}} | {{INS $c.relation}}
$d.Content.Find.Execute(
    "{{c.name}} | {{c.phone}} | {{c.relation}}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{{INS `$c.name}} | {{INS `$c.phone}} | {{INS `$c.relation}}", 2)

# 2) Medications loop body: {{m.name}} | {{m.dosage}} | {{m.schedule}}
#    -> {{INS $m.name}} | {{INS $m.dosage}} | {{INS $m.schedule}}
$d.Content.Find.Execute(
    "{{m.name}} | {{m.dosage}} | {{m.schedule}}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{{INS `$m.name}} | {{INS `$m.dosage}} | {{INS `$m.schedule}}", 2)

# 3) Diagnosis paragraph loop body: {{p}} -> {{INS $p}}
$d.Content.Find.Execute(
    "{{p}}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{{INS `$p}}", 2)

# 4) Each {{END-FOR}} needs its loop variable appended: c, m, then p (in document order).
#    Replace them individually via each paragraph's own range so the three stay distinct.
$loopVars = @("c", "m", "p")
$varIndex = 0
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "{{END-FOR}}") {
        $varName = $loopVars[$varIndex]
        $para.Range.Find.Execute(
            "{{END-FOR}}", $true, $false, $false, $false, $false,
            $true, 1, $false, "{{END-FOR $varName}}", 2)
        $varIndex++
    }
}
